$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the raw poll numbers (row 2) that feed the normalised percentages
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 9
$ws.Range("I2").Value = 9

# Update the raw vote count that feeds the two-party preferred split
$ws.Range("A10").Value = 49

# Leave the selection on H12, matching the last active cell after the edits
$ws.Range("H12").Select()
